# Auto-generated edit script updating the cryptos price/volume table
# to match the commit "Updated cryptos list on Sat Aug 31 06:44:40 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.245.59"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.525.87"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "2.525.06"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "2.975.43"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "59.188.42"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "2.527.75"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.61%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.168"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").Value = "0.0₃0777"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  +6.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "164.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "288.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.26%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "132.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.25%  "
